$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 2"
$ws.Cells.Item($row, 3).Value = "Esportes"
$ws.Cells.Item($row, 4).Value = "2025-04-01T19:38"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Mudança na presidência do Americano. Tolentino Reis é destituído por votação do Conselho Deliberativo. Laila Póvoa assume. *nota coberta*"
